# Append 7 new match rows (rows 9-15) to the "Eoin Morgan (c)" sheet,
# duplicating existing match records (in source order 3,4,5,8,2,6,7),
# and extend the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$batsman = "Eoin Morgan" + $nbsp + "(c)"

$newRows = @(
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", $batsman, "40", "25", "5", "2", "160.00"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", $batsman, "34", "23", "3", "1", "147.82"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", $batsman, "15", "12", "2", "0", "125.00"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", $batsman, "17", "9", "2", "1", "188.88"),
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", $batsman, "39", "29", "2", "2", "134.48"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", $batsman, "68", "35", "5", "6", "194.28"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", $batsman, "30", "34", "3", "1", "88.23")
)

$startRow = 9
$endRow = $startRow + $newRows.Count - 1

# Columns G:K (totalRuns, totalBalls, total4s, total6s, sr) hold numbers
# that must be stored as text (matching the existing rows' t="str" cells),
# so force text formatting on that block before writing the values.
$numRng = $ws.Range("G" + $startRow + ":K" + $endRow)
$numRng.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
